# CryCompanywiseStockReport -- stock-count correction pass
#
# The source stock report stored every figure as a plain literal (no formulas:
# column G = column F (qty) * column D (rate), and each company's "Sub Total:"
# row in column B = SUM(G) over that company's item rows). This script re-enters
# the corrected quantities/values exactly as the corrected report has them, one
# worksheet row at a time, including the knock-on Sub Total / Grand Total rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30: corrected qty (F) / value (G) after stock recount
$ws.Range("F30").Value = 194
$ws.Range("G30").Value = 4979.98

# Row 31: corrected qty (F) / value (G) after stock recount
$ws.Range("F31").Value = 83
$ws.Range("G31").Value = 2592.09

# Row 56: roll up corrected Sub Total / Grand Total
$ws.Range("B56").Value = 47180.85

# Row 104: corrected qty (F) / value (G) after stock recount
$ws.Range("F104").Value = 155
$ws.Range("G104").Value = 15881.3

# Row 106: corrected qty (F) / value (G) after stock recount
$ws.Range("F106").Value = 35
$ws.Range("G106").Value = 4716.95

# Row 114: corrected qty (F) / value (G) after stock recount
$ws.Range("F114").Value = 298
$ws.Range("G114").Value = 5799.08

# Row 115: roll up corrected Sub Total / Grand Total
$ws.Range("B115").Value = 262533.62

# Row 126: corrected qty (F) / value (G) after stock recount
$ws.Range("F126").Value = 45
$ws.Range("G126").Value = 6364.35

# Row 132: corrected qty (F) / value (G) after stock recount
$ws.Range("F132").Value = 1
$ws.Range("G132").Value = 19.28

# Row 133: roll up corrected Sub Total / Grand Total
$ws.Range("B133").Value = 13382.22

# Row 170: corrected qty (F) / value (G) after stock recount
$ws.Range("F170").Value = 2
$ws.Range("G170").Value = 168.84

# Row 185: roll up corrected Sub Total / Grand Total
$ws.Range("B185").Value = 20030.64

# Row 214: corrected qty (F) / value (G) after stock recount
$ws.Range("F214").Value = 50
$ws.Range("G214").Value = 4241.5

# Row 217: roll up corrected Sub Total / Grand Total
$ws.Range("B217").Value = 11245.13

# Row 226: corrected qty (F) / value (G) after stock recount
$ws.Range("F226").Value = 4
$ws.Range("G226").Value = 1114.32

# Row 228: roll up corrected Sub Total / Grand Total
$ws.Range("B228").Value = 26069.78

# Row 236: corrected qty (F) / value (G) after stock recount
$ws.Range("F236").Value = 0
$ws.Range("G236").Value = 0

# Row 238: roll up corrected Sub Total / Grand Total
$ws.Range("B238").Value = 8599.709999999999

# Row 257: corrected qty (F) / value (G) after stock recount
$ws.Range("F257").Value = 0
$ws.Range("G257").Value = 0

# Row 273: corrected qty (F) / value (G) after stock recount
$ws.Range("F273").Value = 8
$ws.Range("G273").Value = 275.68

# Row 279: roll up corrected Sub Total / Grand Total
$ws.Range("B279").Value = 120569.07

# Rows 283-284: two "HUL-Bru Inst Poly 50g" lines had their item code/rate/MRP/
# qty/value swapped onto the wrong row; put each set of figures on its correct row.
$ws.Range("B283").Value = 57077
$ws.Range("D283").Value = 93.08
$ws.Range("E283").Value = 111.2
$ws.Range("F283").Value = 1
$ws.Range("G283").Value = 93.08

$ws.Range("B284").Value = 61610
$ws.Range("D284").Value = 102.71
$ws.Range("E284").Value = 122.71
$ws.Range("F284").Value = 211
$ws.Range("G284").Value = 21671.81

# Row 287: corrected qty (F) / value (G) after stock recount
$ws.Range("F287").Value = 34
$ws.Range("G287").Value = 17055.42

# Row 293: corrected qty (F) / value (G) after stock recount
$ws.Range("F293").Value = 188
$ws.Range("G293").Value = 21475.24

# Row 295: corrected qty (F) / value (G) after stock recount
$ws.Range("F295").Value = 39
$ws.Range("G295").Value = 7486.44

# Row 329: corrected qty (F) / value (G) after stock recount
$ws.Range("F329").Value = 46
$ws.Range("G329").Value = 5456.98

# Row 330: corrected qty (F) / value (G) after stock recount
$ws.Range("F330").Value = 118
$ws.Range("G330").Value = 6977.34

# Row 345: corrected qty (F) / value (G) after stock recount
$ws.Range("F345").Value = 38
$ws.Range("G345").Value = 2342.32

# Row 349: roll up corrected Sub Total / Grand Total
$ws.Range("B349").Value = 378656.91

# Row 357: corrected qty (F) / value (G) after stock recount
$ws.Range("F357").Value = 9
$ws.Range("G357").Value = 1964.07

# Row 358: roll up corrected Sub Total / Grand Total
$ws.Range("B358").Value = 33226.62

# Row 383: corrected qty (F) / value (G) after stock recount
$ws.Range("F383").Value = 1
$ws.Range("G383").Value = 45.97

# Row 389: roll up corrected Sub Total / Grand Total
$ws.Range("B389").Value = 12000.52

# Rows 396-397: same fix for the two "KUS-Floor Wiper" lines.
$ws.Range("B396").Value = 47097
$ws.Range("D396").Value = 112.28
$ws.Range("E396").Value = 134.16
$ws.Range("F396").Value = 15
$ws.Range("G396").Value = 1684.2

$ws.Range("B397").Value = 58047
$ws.Range("D397").Value = 105.54
$ws.Range("E397").Value = 126.1
$ws.Range("F397").Value = 62
$ws.Range("G397").Value = 6543.48

# Row 402: corrected qty (F) / value (G) after stock recount
$ws.Range("F402").Value = 64
$ws.Range("G402").Value = 6182.4

# Row 407: roll up corrected Sub Total / Grand Total
$ws.Range("B407").Value = 50999.09

# Row 417: corrected qty (F) / value (G) after stock recount
$ws.Range("F417").Value = 71
$ws.Range("G417").Value = 13237.24

# Row 418: corrected qty (F) / value (G) after stock recount
$ws.Range("F418").Value = 93
$ws.Range("G418").Value = 3092.25

# Row 424: roll up corrected Sub Total / Grand Total
$ws.Range("B424").Value = 48945.07

# Row 445: corrected qty (F) / value (G) after stock recount
$ws.Range("F445").Value = 132
$ws.Range("G445").Value = 2568.72

# Row 453: roll up corrected Sub Total / Grand Total
$ws.Range("B453").Value = 108531.42

# Row 455: corrected qty (F) / value (G) after stock recount
$ws.Range("F455").Value = 54
$ws.Range("G455").Value = 1996.38

# Row 460: roll up corrected Sub Total / Grand Total
$ws.Range("B460").Value = 9674.01

# Row 515: corrected qty (F) / value (G) after stock recount
$ws.Range("F515").Value = 24
$ws.Range("G515").Value = 625.92

# Row 516: corrected qty (F) / value (G) after stock recount
$ws.Range("F516").Value = 37
$ws.Range("G516").Value = 1930.66

# Row 525: roll up corrected Sub Total / Grand Total
$ws.Range("B525").Value = 28901.76

# Row 564: corrected qty (F) / value (G) after stock recount
$ws.Range("F564").Value = 12
$ws.Range("G564").Value = 328.8

# Row 566: corrected qty (F) / value (G) after stock recount
$ws.Range("F566").Value = 127
$ws.Range("G566").Value = 3383.28

# Row 571: roll up corrected Sub Total / Grand Total
$ws.Range("B571").Value = 45387.06

# Row 603: corrected qty (F) / value (G) after stock recount
$ws.Range("F603").Value = 13
$ws.Range("G603").Value = 430.43

# Row 608: roll up corrected Sub Total / Grand Total
$ws.Range("B608").Value = 33309.74

# Row 662: corrected qty (F) / value (G) after stock recount
$ws.Range("F662").Value = 12
$ws.Range("G662").Value = 2945.64

# Row 664: roll up corrected Sub Total / Grand Total
$ws.Range("B664").Value = 23086.52

# Row 668: corrected qty (F) / value (G) after stock recount
$ws.Range("F668").Value = 70
$ws.Range("G668").Value = 7807.8

# Row 671: corrected qty (F) / value (G) after stock recount
$ws.Range("F671").Value = 66
$ws.Range("G671").Value = 2735.7

# Row 676: roll up corrected Sub Total / Grand Total
$ws.Range("B676").Value = 46197.38

# Row 690: corrected qty (F) / value (G) after stock recount
$ws.Range("F690").Value = 2
$ws.Range("G690").Value = 73.66

# Row 694: roll up corrected Sub Total / Grand Total
$ws.Range("B694").Value = 24963.1

# Row 699: corrected qty (F) / value (G) after stock recount
$ws.Range("F699").Value = 114
$ws.Range("G699").Value = 4263.6

# Row 702: roll up corrected Sub Total / Grand Total
$ws.Range("B702").Value = 10472.54

# Row 747: corrected qty (F) / value (G) after stock recount
$ws.Range("F747").Value = 1163
$ws.Range("G747").Value = 189696.93

# Row 749: corrected qty (F) / value (G) after stock recount
$ws.Range("F749").Value = 218
$ws.Range("G749").Value = 16816.52

# Row 751: corrected qty (F) / value (G) after stock recount
$ws.Range("F751").Value = 117
$ws.Range("G751").Value = 7897.5

# Row 752: roll up corrected Sub Total / Grand Total
$ws.Range("B752").Value = 222073.97

# Row 753: roll up corrected Sub Total / Grand Total
$ws.Range("B753").Value = 2325388.54

# Row 754: roll up corrected Sub Total / Grand Total
$ws.Range("B754").Value = 2325388.54
